# Apply the "Updated cryptos list" data refresh (prices / 1h volume %, and a few
# re-ranked coin rows where Name/Link/Price/Volume shifted to an adjacent row).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '29.405.71'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  +0.15%  '
# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.849.51'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  +0.22%  '
# Row 4
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  +0.08%  '
# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '240.70'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  +0.21%  '
# Row 6
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  +0.00%  '
# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '1.000'
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  +0.03%  '
# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.07688'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  +2.03%  '
# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.2939'
# Row 10
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  +0.28%  '
# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07749'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  +0.45%  '
# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.841.86'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  -0.44%  '
# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '5.017'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  +0.51%  '
# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.00001092'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  +9.28%  '
# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.6801'
# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '83.60'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  +0.98%  '
# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '2.094.54'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  -7.50%  '
# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '6.152'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  +0.44%  '
# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '29.417.42'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  +0.07%  '
# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '229.27'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  +0.66%  '
# Row 21
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  +0.32%  '
# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '7.444'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  -1.11%  '
# Row 24
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  +0.00%  '
# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '157.41'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  +0.04%  '
# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.1389'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  -0.49%  '
# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '8.363'
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  +0.02%  '
# Row 28
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  +0.29%  '
# Row 29
$ws.Range("B29").Value = 'Toncoin'
$ws.Range("C29").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.313'
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  +4.64%  '
# Row 30
$ws.Range("B30").Value = 'PancakeSwap'
$ws.Range("C30").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.468'
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  +0.06%  '
# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.05680'
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  +0.08%  '
# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.112'
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  -0.29%  '
# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '4.049'
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  +0.84%  '
# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.849'
# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.158'
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  +0.40%  '
# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.7096'
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  -0.50%  '
# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.586'
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  -0.18%  '
# Row 38
$ws.Range("B38").Value = 'Maker'
$ws.Range("C38").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.232.69'
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  -2.00%  '
# Row 39
$ws.Range("B39").Value = 'MXToken'
$ws.Range("C39").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.776'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  -0.19%  '
# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.01799'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  -0.74%  '
# Row 41
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  +4.09%  '
# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.9139'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  +0.08%  '
# Row 44
$ws.Range("B44").Value = 'RocketPoolETH'
$ws.Range("C44").Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '2.003.79'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  +0.09%  '
# Row 45
$ws.Range("B45").Value = 'Quant'
$ws.Range("C45").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '101.48'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  +0.44%  '
# Row 46
$ws.Range("B46").Value = 'Aave'
$ws.Range("C46").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '66.23'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  -0.06%  '
# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '7.169'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  +1.88%  '
# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.4014'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  -0.70%  '
# Row 49
$ws.Range("B49").Value = 'RenderToken'
$ws.Range("C49").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.688'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  +0.15%  '
# Row 50
$ws.Range("B50").Value = 'EnergySwap'
$ws.Range("C50").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '8.996'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  -1.27%  '
# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.1124'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  -0.03%  '
